# MQTT_Topics.xlsx update
#
# Adds a new "Alerts" topic block (columns M:N) mirroring the existing
# Topic/SubTopic header layout, changes the DAQ block's second header
# (G2) from "SubTopic" to "Value", and gives it a new "As JSON" value
# row (G3) styled like the other section headers (italic, on the DAQ
# fill). Selection moves to G3 to match the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Alerts" header block: M2:N2 ("Topic"/"SubTopic"), same look
#     as the other header cells (A2:C2, F2:G2, I2:K2). ---
$ws.Range("A2").Copy()
$ws.Range("M2").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("M2").Value = "Topic"

$ws.Range("B2").Copy()
$ws.Range("N2").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("N2").Value = "SubTopic"

# --- DAQ block: G2 changes from "SubTopic" to "Value" ---
$ws.Range("G2").Value = "Value"

# --- New "Alerts" section row: M3 "Alerts" + N3, on a new red fill ---
$ws.Range("M3:N3").Interior.Color = 5263615  # RGB(255,80,80) -> FFFF5050
$ws.Range("M3").Value = "Alerts"

# --- DAQ block: G3 becomes an italic "As JSON" label on the DAQ header fill ---
$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("G3").Font.Italic = $true
$ws.Range("G3").Value = "As JSON"

# --- New "Alerts" data row: N4 "x_axis" (plain, unstyled like B4/G4/J4) ---
$ws.Range("N4").Value = "x_axis"

# --- Match the saved selection state ---
$null = $ws.Range("G3").Select()
